# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4316
$ws1.Range("F3").Value = 2440
$ws1.Range("F9").Value = 128
$ws1.Range("F10").Value = 141
$ws1.Range("F12").Value = 1604
$ws1.Range("F14").Value = 3374
$ws1.Range("F15").Value = 228

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4316
$ws4.Range("F3").Value = 2440
$ws4.Range("F11").Value = 128
$ws4.Range("F12").Value = 141
$ws4.Range("F16").Value = 1604
$ws4.Range("F18").Value = 3374
$ws4.Range("F19").Value = 228
